# Austria Landesliga base update (04-04-2024 23:22)
# Several rows had their match-data (all columns except the leading
# row-sequence-number column A) swapped with the adjacent row, as the
# underlying data source was re-sorted. Column A (and the row number
# itself) must stay where it is; only columns B:AC trade places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is a pair of worksheet rows whose B:AC content must be
# exchanged.
$rowPairs = @(
    @(465, 466),
    @(485, 486),
    @(491, 492),
    @(496, 497),
    @(509, 510),
    @(518, 519),
    @(531, 532)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1" + ":AC$r1")
    $range2 = $ws.Range("B$r2" + ":AC$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
